# Sync attendance_reports: reorder "Recorded By" (column G) entries so that
# any comma-separated list ending in exactly "System" has its order reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Length -gt 1 -and $parts[$parts.Length - 1] -eq "System") {
            $reversed = $parts[($parts.Length - 1)..0]
            $cell.Value = [string]::Join(", ", $reversed)
        }
    }
}
